# Update the "Metadata" sheet of the ValueSet workbook:
#  - Version bump 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date refreshed to the new publication timestamp
#  - Contact rows replaced with the CIBMTR org contact + Bob Milius
#  - A new "Jurisdiction" row inserted after the Contact rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version (row 3)
$ws.Range("B3").Value = "0.1.7"

# Status (row 6)
$ws.Range("B6").Value = "draft"

# Date (row 8)
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Contact (row 10) - was "No display for ContactDetail"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Contact (row 11) - second contact, was "No display for ContactDetail"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new "Jurisdiction" row (row 12), pushing Description/Purpose/Copyright/Immutable down
$ws.Rows.Item(12).Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$excel.CutCopyMode = $false
